$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 823 ("「がんばって！」بالتوفيق ..." post), shifting all subsequent
# rows up by one. Excel will automatically update the sheet's used range /
# dimension accordingly.
$ws.Rows.Item(823).Delete()
